$d = $word.ActiveDocument

function Set-ParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}

# 1) "closekb" -> explanatory sentence
Set-ParaText 3 "Crear una función para cerrar el teclado."

# 2) "readchar" -> explanatory sentence
Set-ParaText 4 "Crear una función para leer un carácter."

# 3) "readString" -> explanatory sentence
Set-ParaText 5 "Crear una función para leer una cadena."

# 4) readBoolean (3 args) -> explanatory sentence (note: "funcion" keeps the
#    author's original (unaccented) typo, mirroring the source edit)
Set-ParaText 6 "Crear una funcion que devuelva true si se introduce 1, o false si se introduce 2. Enbuclar si hay un error."

# 5) readBoolean (1 arg) -> explanatory sentence (leading space preserved)
Set-ParaText 7 " Crear función que devuelva true si devuelva true si se introduce s o S, o false si se introduce n o N. Enbuclar si hay un error."

# 6) leer numero -> explanatory sentence
Set-ParaText 8 "Crear función que devuelva un número introducido por teclado (todos los tipos). Enbuclar si hay un error."

# 7) leer numero con equivalencias -> explanatory sentence (leading space preserved)
Set-ParaText 9 " Crear función para leer número con equivalencias (todos los tipos). Enbuclar si hay un error."

# 8) leer numero en rango -> append "Enbuclar si hay un error." and move the
#    "_GoBack" bookmark here (it previously sat at the end of the
#    "Incluir InputMismatchException" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$p10 = $d.Paragraphs.Item(10)
$r10 = $p10.Range
$r10.End = $r10.End - 1
$r10.InsertAfter(" Enbuclar si hay un error.")
$bmRange = $p10.Range.Duplicate
$bmRange.End = $bmRange.End - 1
$d.Bookmarks.Add("_GoBack", $bmRange)

# 9) "Incluir InputMismatchException al introducir numero." -> updated wording
Set-ParaText 13 "Incluir InputMismatchException al introducir dato erróneo."

# 10) "Pruebas" -> "Incluir Excepciones si las hay."
Set-ParaText 15 "Incluir Excepciones si las hay."
